$wb = $excel.ActiveWorkbook

# --- Update status text "Handed back: in sync with en-US" -> "Ready for handoff" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Update the "Latest HO Xliff Generate Date" / handback timestamp 19:01:32 -> 19:02:14 ---
$wsOverview.Range("G2").Value = "2016-08-21 19:02:14"
$wsDeDe.Range("H2").Value = "2016-08-21 19:02:14"

# --- Update the "Latest Handoff Datetime" for zh-cn 19:01:28 -> 19:02:10 ---
$wsZhCn.Range("H2").Value = "2016-08-21 19:02:10"

# --- Update column widths ---
# Overview sheet columns E and F: 29.9777047293527 -> 17.2159881591797
# (COM ColumnWidth quantizes to whole pixels; 16.26 is the nearest input
# that rounds to the target stored width.)
$wsOverview.Columns.Item(5).ColumnWidth = 16.26
$wsOverview.Columns.Item(6).ColumnWidth = 16.26

# zh-cn sheet column C: 29.9777047293527 -> 17.2159881591797
$wsZhCn.Columns.Item(3).ColumnWidth = 16.26

# de-de sheet column C: 29.9777047293527 -> 17.2159881591797
$wsDeDe.Columns.Item(3).ColumnWidth = 16.26
